$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (SOC, Initial SOC) before the existing Dc column (C)
$ws.Columns("C:D").Insert()

# Insert one new column (Dt*) after the Dc column, which is now column E
$ws.Columns("F:F").Insert()

# Header row
$ws.Range("C1").Value = "SOC"
$ws.Range("D1").Value = "Initial SOC"
$ws.Range("F1").Value = "Dt* (cm^2/s)"

$ws.Range("C2").Value = [double]"0.9916058068842946"
$ws.Range("D2").Value = [double]"0.99161342784741"
$ws.Range("C3").Value = [double]"0.9883656497179512"
$ws.Range("D3").Value = [double]"0.9915982035525202"
$ws.Range("C4").Value = [double]"0.9599077284171872"
$ws.Range("D4").Value = [double]"0.9851331141320233"
$ws.Range("F4").Value = [double]"1.15229688516344e-12"
$ws.Range("C5").Value = [double]"0.9146849503939857"
$ws.Range("D5").Value = [double]"0.9346823588829911"
$ws.Range("F5").Value = [double]"1.776715428762427e-13"
$ws.Range("C6").Value = [double]"0.8744998873061952"
$ws.Range("D6").Value = [double]"0.8946875594981116"
$ws.Range("F6").Value = [double]"1.08987751737354e-13"
$ws.Range("C7").Value = [double]"0.8280769967696261"
$ws.Range("D7").Value = [double]"0.8543122315189713"
$ws.Range("F7").Value = [double]"4.194112058806839e-13"
$ws.Range("C8").Value = [double]"0.7844593898715593"
$ws.Range("D8").Value = [double]"0.8018417793085657"
$ws.Range("F8").Value = [double]"3.243662597313021e-13"
$ws.Range("C9").Value = [double]"0.753023311525908"
$ws.Range("D9").Value = [double]"0.7670770671534796"
$ws.Range("F9").Value = [double]"4.464375925515226e-13"
$ws.Range("C10").Value = [double]"0.7236103849186528"
$ws.Range("D10").Value = [double]"0.7389695719200885"
$ws.Range("F10").Value = [double]"4.479915614791097e-13"
$ws.Range("C11").Value = [double]"0.6922217245824136"
$ws.Range("D11").Value = [double]"0.7082512141137918"
$ws.Range("F11").Value = [double]"4.797775399747223e-13"
$ws.Range("C12").Value = [double]"0.6608854425592587"
$ws.Range("D12").Value = [double]"0.6761922513716577"
$ws.Range("F12").Value = [double]"5.998924816421758e-13"
$ws.Range("C13").Value = [double]"0.6325963643331218"
$ws.Range("D13").Value = [double]"0.6455786498829891"
$ws.Range("F13").Value = [double]"7.170888377214106e-13"
$ws.Range("C14").Value = [double]"0.6091086802074966"
$ws.Range("D14").Value = [double]"0.61961409468437"
$ws.Range("F14").Value = [double]"7.418080744046689e-13"
$ws.Range("C15").Value = [double]"0.5896660542308286"
$ws.Range("D15").Value = [double]"0.5986032820345555"
$ws.Range("F15").Value = [double]"7.083580255877923e-13"
$ws.Range("C16").Value = [double]"0.5723526510347539"
$ws.Range("D16").Value = [double]"0.5807288436250796"
$ws.Range("F16").Value = [double]"6.121501612891268e-13"
$ws.Range("C17").Value = [double]"0.5551252160970644"
$ws.Range("D17").Value = [double]"0.5639764754164138"
$ws.Range("F17").Value = [double]"5.033052173925559e-13"
$ws.Range("C18").Value = [double]"0.5356900910255624"
$ws.Range("D18").Value = [double]"0.5462739734253723"
$ws.Range("F18").Value = [double]"4.23244414795789e-13"
$ws.Range("C19").Value = [double]"0.511780929470421"
$ws.Range("D19").Value = [double]"0.525106225722237"
$ws.Range("F19").Value = [double]"3.426131966518981e-13"
$ws.Range("C20").Value = [double]"0.4649917654757544"
$ws.Range("D20").Value = [double]"0.4984556499488101"
$ws.Range("F20").Value = [double]"3.593866031800972e-13"
$ws.Range("C21").Value = [double]"0.4202710072270374"
$ws.Range("D21").Value = [double]"0.4315278963900309"
$ws.Range("F21").Value = [double]"5.268185567886848e-13"
$ws.Range("C22").Value = [double]"0.4009009135705802"
$ws.Range("D22").Value = [double]"0.4090141332682362"
$ws.Range("F22").Value = [double]"4.292087019014612e-13"
$ws.Range("C23").Value = [double]"0.3862656649152406"
$ws.Range("D23").Value = [double]"0.3927877089106672"
$ws.Range("F23").Value = [double]"3.053995844799529e-13"
$ws.Range("C24").Value = [double]"0.3738045745778091"
$ws.Range("D24").Value = [double]"0.3797436357360754"
$ws.Range("F24").Value = [double]"1.896999292553254e-13"
$ws.Range("C25").Value = [double]"0.3615206605362039"
$ws.Range("D25").Value = [double]"0.367865528517731"
$ws.Range("F25").Value = [double]"7.515273775318008e-14"
$ws.Range("C26").Value = [double]"0.3478486156914935"
$ws.Range("D26").Value = [double]"0.355175807637528"
$ws.Range("F26").Value = [double]"8.487107204345371e-14"
$ws.Range("C27").Value = [double]"0.3121528774689605"
$ws.Range("D27").Value = [double]"0.3405214395329013"
